$wb = $excel.ActiveWorkbook
$wsProduct = $wb.Worksheets.Item("Product Backlog")
$wsSprint = $wb.Worksheets.Item("Sprint Backlog #1")

# --- Product Backlog: double the "Pracochlonnosc" (effort) values ---
$wsProduct.Range("C2").Value = 12
$wsProduct.Range("C3").Value = 8
$wsProduct.Range("C4").Value = 6
$wsProduct.Range("C5").Value = 2
$wsProduct.Range("C6").Value = 4
$wsProduct.Range("C7").Value = 6

# New summary row with a SUM formula
$wsProduct.Range("C10").Formula = "=SUM(C2:C4)"

# --- Sprint Backlog #1: add a new "oby poszlo" note cell ---
$wsSprint.Range("B2").Value = "oby poszło"

# Add a new "Kolumna1" column to the Tabela3 table
$loSprint = $wsSprint.ListObjects.Item("Tabela3")
$loSprint.ListColumns.Add()
$wsSprint.Range("D4").Value = "Kolumna1"
$wsProduct.Range("F1").Copy()
$wsSprint.Range("D4").PasteSpecial(-4122)

# Double the "Pracochlonnosc" (effort) values, matching the Product Backlog change
$wsSprint.Range("C6").Value = 3
$wsSprint.Range("C7").Value = 3
$wsSprint.Range("C8").Value = 6
$wsSprint.Range("C9").Value = 0
$wsSprint.Range("C10").Value = 2
$wsSprint.Range("C11").Value = 2
$wsSprint.Range("C12").Value = 2
$wsSprint.Range("C13").Value = 2
$wsSprint.Range("C14").Value = 0
$wsSprint.Range("C15").Value = 2
$wsSprint.Range("C16").Value = 2
$wsSprint.Range("C17").Value = 1
$wsSprint.Range("C18").Value = 1

# --- Update selections to match the new focus area ---
$wsProduct.Range("B12").Select()
$wsSprint.Range("D6").Select()

# Sprint Backlog #1 becomes the active/selected sheet
$wsSprint.Activate()
